$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at D:E (old E..H shift right to G..J)
$ws.Range("D1:E1").EntireColumn.Insert()

# Header row
$ws.Range("D1").Value = "Date Created"
$ws.Range("E1").Value = "Date Completed"

# Row 2 (Issue #1) - default formatting, gets a fresh date style
$ws.Range("D2").Value = 43542
$ws.Range("D2").NumberFormat = "mm-dd-yy"

# Row 3 (Issue #2) - "Neutral" row style; format D3 then clone format onto E3
$ws.Range("D3").Value = 43542
$ws.Range("D3").NumberFormat = "mm-dd-yy"
$ws.Range("E3").Value = 43542
$ws.Range("D3").Copy()
$ws.Range("E3").PasteSpecial(-4122)

# Row 4 (Issue #3) - "Good" row style; format D4 then clone format onto E4
$ws.Range("D4").Value = 43542
$ws.Range("D4").NumberFormat = "mm-dd-yy"
$ws.Range("E4").Value = 43542
$ws.Range("D4").Copy()
$ws.Range("E4").PasteSpecial(-4122)

# Rows 5-14 (Issues #4-#13) - plain values, no explicit number format
$ws.Range("D5").Value = 43542
$ws.Range("D6").Value = 43542
$ws.Range("D7").Value = 43542
$ws.Range("D8").Value = 43542
$ws.Range("D9").Value = 43542
$ws.Range("D10").Value = 43542
$ws.Range("D11").Value = 43542
$ws.Range("D12").Value = 43542
$ws.Range("D13").Value = 43542
$ws.Range("D14").Value = 43542

# Row 15 (Issue #14) - "Good" row style again; clone format from D4 (same resulting style)
$ws.Range("D15").Value = 43542
$ws.Range("D4").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D15").Value = 43542
$ws.Range("E15").Value = 43542
$ws.Range("D15").Copy()
$ws.Range("E15").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Column widths for the newly inserted columns
$ws.Columns("D").ColumnWidth = 12.5703125
$ws.Columns("E").ColumnWidth = 15.5703125

# Match the saved selection from the authored edit
$ws.Range("E20").Select()
